# Proof-reading pass: promote a few "second_pass" rows to "proofread"
# (plateK10rep1 WellA06, plateK10rep1 WellG06, plateK11rep1 WellD05),
# then leave the selection on the last-edited cell (B38) as the user
# would after working their way down the Status column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B35").Value = "proofread"
$ws.Range("B36").Value = "proofread"
$ws.Range("B38").Value = "proofread"

$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("B38").Select()
